$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (shared between input and output sheets)
$wsInput.Range("B1").Value = "4136-RBI-EPP-DB-SAR-NOREC-MOREREPAY-1st"
$wsOutput.Range("B1").Value = "4136-RBI-EPP-DB-SAR-NOREC-MOREREPAY-1st"

# Update short name - was numeric 4136, now text "413z"
$wsInput.Range("B2").Value = "413z"

# Select B1 on input sheet (was B18)
$wsInput.Range("B1").Select() | Out-Null

# Make the output sheet the active sheet/tab
$wsOutput.Activate() | Out-Null
$wsOutput.Range("B1").Select() | Out-Null
